$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Ingresar al SIU " -> "Ingresa" + " al SIU " (fix "Ingresar" typo)
# ---------------------------------------------------------------------------

# Find the paragraph that currently starts with "Ingresar al SIU "
$target1 = $d.Content
$found1 = $target1.Find.Execute("Ingresar al SIU ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $insertPoint1 = $d.Range($target1.Start, $target1.Start)

    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Ingresa</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> al SIU </w:t></w:r><w:r w:rsidR="00333EB1"><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Guaran&#237;</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">  de grado con usuario y contrase&#241;a </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint1.InsertXML($xml1)

    # Remove the old (un-split, mis-spelled) run that now follows what we just inserted
    $p1 = $insertPoint1.Paragraphs(1)
    $oldR1 = $d.Range($p1.Range.Start, $p1.Range.End)
    $oldR1.Find.Execute("Ingresar al SIU Guaran" + [char]0x00ED + "  de grado con usuario y contrase" + [char]0x00F1 + "a ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $oldR1.Text = ""
}

# ---------------------------------------------------------------------------
# Change 2: remove the stray "de " before "para realizarlo" and split the run
# "Estate atento ... tiene un periodo de para realizarlo"
#   -> "Estate atento ... cada c" + "uatrimestre tiene un periodo " + "para realizarlo"
# ---------------------------------------------------------------------------

$target2 = $d.Content
$found2 = $target2.Find.Execute("tiene un periodo de para realizarlo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $p2 = $target2.Paragraphs(1)
    $insertPoint2 = $d.Range($p2.Range.Start, $p2.Range.Start)

    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Estate atento a las fechas de inscripci&#243;n, ya que cada c</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">uatrimestre tiene un periodo </w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>para realizarlo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint2.InsertXML($xml2)

    # Remove the old (un-split) run that now follows what we just inserted
    $p2b = $insertPoint2.Paragraphs(1)
    $oldR2 = $d.Range($p2b.Range.Start, $p2b.Range.End)
    $oldR2.Find.Execute("Estate atento a las fechas de inscripci" + [char]0x00F3 + "n, ya que cada cuatrimestre tiene un periodo de para realizarlo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $oldR2.Text = ""
}

Write-Output "done"
